# Update countries & provincias Spain
# Applies the 2020-04-15 17:22 -> 17:52 COVID data refresh to the "Pais" sheet:
#  - updates case counts for several countries
#  - Moldavia overtakes Marruecos (table is sorted desc by "Casos totales")
#  - Republica de Chipre overtakes Crucero
#  - refreshes the "Datos actualizados..." timestamp string in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header ---
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 17:52"

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 616458
$ws.Range("C4").Value = 2572
$ws.Range("D4").Value = 40086
$ws.Range("E4").Value = 550159
$ws.Range("F4").Value = 13473
$ws.Range("G4").Value = 166
$ws.Range("H4").Value = 26213

# --- Alemania (row 8) ---
$ws.Range("B8").Value = 132718
$ws.Range("C8").Value = 508
$ws.Range("D8").Value = 72600
$ws.Range("E8").Value = 56526
$ws.Range("F8").Value = 4288

# --- Canada (row 15) ---
$ws.Range("B15").Value = 27557
$ws.Range("C15").Value = 494
$ws.Range("D15").Value = 8235
$ws.Range("E15").Value = 18368
$ws.Range("F15").Value = 557
$ws.Range("G15").Value = 51
$ws.Range("H15").Value = 954

# --- Polonia (row 30) ---
$ws.Range("B30").Value = 7582
$ws.Range("C30").Value = 380
$ws.Range("D30").Value = 668
$ws.Range("E30").Value = 6628
$ws.Range("F30").Value = 160
$ws.Range("G30").Value = 23
$ws.Range("H30").Value = 286

# --- Chequia (row 35) ---
$ws.Range("B35").Value = 6216
$ws.Range("C35").Value = 105
$ws.Range("D35").Value = 819
$ws.Range("E35").Value = 5231
$ws.Range("F35").Value = 84
$ws.Range("G35").Value = 5
$ws.Range("H35").Value = 166

# --- Moldavia now ranks above Marruecos (row 59 becomes Moldavia, row 60 becomes Marruecos) ---
$ws.Range("A59").Value = "Moldavia"
$ws.Range("B59").Value = 2049
$ws.Range("C59").Value = 115
$ws.Range("D59").Value = 171
$ws.Range("E59").Value = 1834
$ws.Range("F59").Value = 80
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 44

$ws.Range("A60").Value = "Marruecos"
$ws.Range("B60").Value = 1988
$ws.Range("C60").Value = 100
$ws.Range("D60").Value = 218
$ws.Range("E60").Value = 1643
$ws.Range("F60").Value = 1
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 127

# --- Lituania (row 76) ---
$ws.Range("E76").Value = 923
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 30

# --- Cuba (row 82) ---
$ws.Range("F82").Value = 15

# --- Republica de Chipre now ranks above Crucero (row 86 becomes Chipre, row 87 becomes Crucero) ---
$ws.Range("A86").Value = "Republica de Chipre"
$ws.Range("B86").Value = 715
$ws.Range("C86").Value = 20
$ws.Range("D86").Value = 65
$ws.Range("E86").Value = 638
$ws.Range("F86").Value = 8

$ws.Range("A87").Value = "Crucero"
$ws.Range("B87").Value = 712
$ws.Range("D87").Value = 639
$ws.Range("E87").Value = 61
$ws.Range("F87").Value = 7

# --- Congo (row 132) ---
$ws.Range("D132").Value = 11
$ws.Range("E132").Value = 101

# --- Tanzania (row 138) ---
$ws.Range("D138").Value = 11
$ws.Range("E138").Value = 73
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 4
